$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.606.18"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.423.14"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").Value = "3.416.46"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000282"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "694.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "3.973.23"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "69.646.77"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "3.425.78"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.896"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "572.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "3.558.52"
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "0.0₃0732"
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0418"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.82%  "
